# Applies the value updates to the three-digit x one-digit multiplication
# table. Each populated cell contains a single text run like "959×8=7672".
# We target each cell by its (row, column) coordinates -- which match the
# order the values appear in the document -- and overwrite the cell's text
# directly (rather than using Find/Replace, since some values such as
# "234×5=1170" repeat verbatim in two different cells and this runtime's
# Find/Replace updates every run sharing identical text instead of just the
# one in scope). Writing directly to a Range limited to the cell's content
# (excluding the trailing cell-mark) preserves the existing run formatting
# (font/size) because it replaces the text of the existing run in place.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# List of (row, column, oldText, newText) tuples in document order.
$updates = @(
    @(1, 1, "959×8=7672", "848×6=5088"),
    @(1, 2, "234×5=1170", "157×3=471"),
    @(1, 3, "301×2=602", "155×3=465"),
    @(1, 4, "809×8=6472", "204×6=1224"),
    @(1, 5, "795×6=4770", "611×4=2444"),

    @(5, 1, "234×5=1170", "682×3=2046"),
    @(5, 2, "580×4=2320", "483×3=1449"),
    @(5, 3, "280×2=560", "478×4=1912"),
    @(5, 4, "260×9=2340", "234×9=2106"),
    @(5, 5, "931×7=6517", "746×3=2238"),

    @(10, 1, "944×4=3776", "443×3=1329"),
    @(10, 2, "569×4=2276", "106×2=212"),
    @(10, 3, "684×8=5472", "548×6=3288"),
    @(10, 4, "761×2=1522", "301×6=1806"),
    @(10, 5, "406×5=2030", "424×2=848"),

    @(15, 1, "727×4=2908", "445×8=3560"),
    @(15, 2, "186×6=1116", "370×2=740"),
    @(15, 3, "764×2=1528", "189×5=945"),
    @(15, 4, "687×9=6183", "906×7=6342"),
    @(15, 5, "548×2=1096", "762×8=6096"),

    @(20, 1, "870×4=3480", "387×7=2709"),
    @(20, 2, "108×7=756", "843×3=2529"),
    @(20, 3, "473×3=1419", "795×4=3180"),
    @(20, 4, "877×8=7016", "158×7=1106"),
    @(20, 5, "513×7=3591", "665×3=1995")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $old = $u[2]
    $new = $u[3]

    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the final cell-mark character so we only replace the
    # visible text content, keeping it inside the existing run/paragraph.
    $rng.End = $rng.End - 1

    if ($rng.Text -ne $old) {
        Write-Host ("MISMATCH row=" + $row + " col=" + $col + `
                     " expected=" + $old + " actual=" + $rng.Text)
    }

    $rng.Text = $new
}

Write-Host "Done."
